$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold, centered, bordered) from existing header cell H1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for columns I (I0) and J (IF), rows 2-10
$values = @(
    @(6, 7),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(7, 7),
    @(9, 9),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
